$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.115.87"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.824.75"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'311.73"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.4639"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").Value = "'0.3636"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "'0.07295"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "'0.8700"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").Value = "'20.17"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.860.53"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'0.07623"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "'5.348"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "'92.71"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "'6.482"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'0.000008649"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "27.306.03"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "'5.192"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "2.082.97"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'151.80"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'1.863"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").Value = "'18.27"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'2.107"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").Value = "'116.22"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'5.088"
$ws.Range("E30").Value = "  -3.66%  "
$ws.Range("D31").Value = "'0.08918"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'2.960"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'0.7346"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.142"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.445"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "'2.535"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("D38").Value = "'0.05266"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'1.067"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("D41").Value = "'2.933"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'7.145"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").Value = "'0.5212"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").Value = "'8.271"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("D46").Value = "'0.4880"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'103.96"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'10.14"
$ws.Range("E49").Value = "  -4.06%  "
$ws.Range("D50").Value = "'1.637"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  -1.32%  "
